$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D30").Value = "$473.23"
$ws.Range("F30").Value = "PASS"
$ws.Range("D31").Value = "$252.98"
$ws.Range("F31").Value = "PASS"
